# Rename 'variable' -> 'variable-code' and 'en_long_name'/'en_code_label'
# to 'en_variable-label'/'en_code-label'.
# Close #144

$wb = $excel.ActiveWorkbook

$wsVariables = $wb.Worksheets.Item("Variables")
$wsCodelists = $wb.Worksheets.Item("Codelists")

# --- Content renames ---------------------------------------------------
$wsVariables.Range("C1").Value = "variable-code"
$wsVariables.Range("D1").Value = "en_variable-label"

$wsCodelists.Range("A1").Value = "variable-code"
$wsCodelists.Range("D1").Value = "en_code-label"

# --- Column width adjustments (headers got longer) ----------------------
$wsVariables.Columns.Item(3).ColumnWidth = 15.666666666666666
$wsVariables.Columns.Item(4).ColumnWidth = 18.166666666666668
$wsCodelists.Columns.Item(3).ColumnWidth = 6.333333333333333

# --- Selection / active sheet changes -----------------------------------
$wsVariables.Activate()
$wsVariables.Range("D1").Select()

$wsCodelists.Activate()
$wsCodelists.Range("A2").Select()
